# Apply scheduled-runner market data updates to the Gungnir Profits workbook.
# Each block updates the currentAveragePrice / Leve price / profit columns (H:N)
# for one leve row on one job sheet, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 286.30435
$ws.Range("I39").Value = 53.35294
$ws.Range("J39").Value = 946.3333
$ws.Range("K39").Value = 160.05882
$ws.Range("L39").Value = 2838.9999
$ws.Range("M39").Value = 135.94118
$ws.Range("N39").Value = -3430.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1531.3334
$ws.Range("I127").Value = 658.8
$ws.Range("J127").Value = 2154.5715
$ws.Range("K127").Value = 1976.4
$ws.Range("L127").Value = 6463.7145
$ws.Range("M127").Value = 2983.6
$ws.Range("N127").Value = -16383.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 186.85715
$ws.Range("I5").Value = 41
$ws.Range("J5").Value = 211.16667
$ws.Range("K5").Value = 41
$ws.Range("L5").Value = 211.16667
$ws.Range("M5").Value = 71
$ws.Range("N5").Value = -435.16667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1323.65
$ws.Range("I122").Value = 1123.2667
$ws.Range("J122").Value = 1924.8
$ws.Range("K122").Value = 3369.800099999999
$ws.Range("L122").Value = 5774.4
$ws.Range("M122").Value = -919.8000999999995
$ws.Range("N122").Value = -10674.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 186.85715
$ws.Range("I4").Value = 41
$ws.Range("J4").Value = 211.16667
$ws.Range("K4").Value = 41
$ws.Range("L4").Value = 211.16667
$ws.Range("M4").Value = 74
$ws.Range("N4").Value = -441.16667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 71430290
$ws.Range("I105").Value = 1572.7142
$ws.Range("J105").Value = 142859000
$ws.Range("K105").Value = 1572.7142
$ws.Range("L105").Value = 142859000
$ws.Range("M105").Value = 174.2858000000001
$ws.Range("N105").Value = -142862494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 43480420
$ws.Range("I99").Value = 90910850
$ws.Range("J99").Value = 2530.25
$ws.Range("K99").Value = 90910850
$ws.Range("L99").Value = 2530.25
$ws.Range("M99").Value = -90909352
$ws.Range("N99").Value = -5526.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 648.3333
$ws.Range("I122").Value = 513.3333
$ws.Range("J122").Value = 783.3333
$ws.Range("K122").Value = 1539.9999
$ws.Range("L122").Value = 2349.9999
$ws.Range("M122").Value = 910.0001
$ws.Range("N122").Value = -7249.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 43480420
$ws.Range("I126").Value = 90910850
$ws.Range("J126").Value = 2530.25
$ws.Range("K126").Value = 272732550
$ws.Range("L126").Value = 7590.75
$ws.Range("M126").Value = -272730080
$ws.Range("N126").Value = -12530.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15875854
$ws.Range("I5").Value = 24691684
$ws.Range("J5").Value = 7360.3335
$ws.Range("K5").Value = 74075052
$ws.Range("L5").Value = 22081.0005
$ws.Range("M5").Value = -74074940
$ws.Range("N5").Value = -22305.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27047.297
$ws.Range("J12").Value = 38489.117
$ws.Range("L12").Value = 115467.351
$ws.Range("N12").Value = -115813.351

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1031
$ws.Range("I70").Value = 708
$ws.Range("K70").Value = 2124
$ws.Range("M70").Value = -1809

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1031
$ws.Range("I73").Value = 708
$ws.Range("K73").Value = 2124
$ws.Range("M73").Value = -1032

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 1038.4
$ws.Range("I82").Value = 564
$ws.Range("J82").Value = 1750
$ws.Range("K82").Value = 1692
$ws.Range("L82").Value = 5250
$ws.Range("M82").Value = -1286
$ws.Range("N82").Value = -6062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 1038.4
$ws.Range("I85").Value = 564
$ws.Range("J85").Value = 1750
$ws.Range("K85").Value = 1692
$ws.Range("L85").Value = 5250
$ws.Range("M85").Value = -288
$ws.Range("N85").Value = -8058

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 2000
$ws.Range("J88").Value = 2000
$ws.Range("L88").Value = 6000
$ws.Range("N88").Value = -6856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 2000
$ws.Range("J91").Value = 2000
$ws.Range("L91").Value = 6000
$ws.Range("N91").Value = -8964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 1408
$ws.Range("J94").Value = 1500
$ws.Range("L94").Value = 4500
$ws.Range("N94").Value = -5852

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 922.3333
$ws.Range("I98").Value = 701
$ws.Range("J98").Value = 977.6667
$ws.Range("K98").Value = 2103
$ws.Range("L98").Value = 2933.0001
$ws.Range("M98").Value = -605
$ws.Range("N98").Value = -5929.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 250000940
$ws.Range("I116").Value = 1628.5
$ws.Range("J116").Value = 500000260
$ws.Range("K116").Value = 4885.5
$ws.Range("L116").Value = 1500000780
$ws.Range("M116").Value = -1443.5
$ws.Range("N116").Value = -1500007664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 14874.875
$ws.Range("I118").Value = 981.3333
$ws.Range("J118").Value = 56555.5
$ws.Range("K118").Value = 2943.9999
$ws.Range("L118").Value = 169666.5
$ws.Range("M118").Value = -1700.9999
$ws.Range("N118").Value = -172152.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 12759804
$ws.Range("I122").Value = 48077160
$ws.Range("J122").Value = 6315.222
$ws.Range("K122").Value = 432694440
$ws.Range("L122").Value = 56836.998
$ws.Range("M122").Value = -432691990
$ws.Range("N122").Value = -61736.998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1341.0555
$ws.Range("I129").Value = 1356.5
$ws.Range("J129").Value = 1333.3334
$ws.Range("K129").Value = 4069.5
$ws.Range("L129").Value = 4000.0002
$ws.Range("M129").Value = 930.5
$ws.Range("N129").Value = -14000.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 15875854
$ws.Range("I135").Value = 24691684
$ws.Range("J135").Value = 7360.3335
$ws.Range("K135").Value = 222225156
$ws.Range("L135").Value = 66243.0015
$ws.Range("M135").Value = -222222621
$ws.Range("N135").Value = -71313.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 10000
$ws.Range("J118").Value = 10000
$ws.Range("L118").Value = 10000
$ws.Range("N118").Value = -13314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 26322022
$ws.Range("I122").Value = 62513030
$ws.Range("K122").Value = 187539090
$ws.Range("M122").Value = -187536640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 111111610
$ws.Range("I113").Value = 142857460
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 428572380
$ws.Range("L113").Value = 3450
$ws.Range("M113").Value = -428570210
$ws.Range("N113").Value = -7790

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 21895.04
$ws.Range("I122").Value = 27440.842
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 82322.526
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -79872.526
$ws.Range("N122").Value = -17900.0005
